$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(15, 1, 28.56, 427.516, 2.02334066, 22547074.256, 8.811091039999999, 16.32474726, 466.49302304, 4.917190001924886, 73.4991836962018, 0.3356641017925939, 9958242.651426483, 4.189067863232973, 7.21133636173174, 224.6857040280042, 0.15),
    @(15, 1, 42.998, 629.062, 1.36893166, 27067272.296, 5.544071260000001, 15.57071668, 673.4094539999999, 7.076591856579636, 99.21435951649609, 0.2077273846977117, 9260655.074927676, 2.435910766544815, 6.516669463905607, 316.8585653648009, 0.85),
    @(20, 1, 20.51, 409.236, 2.1201033, 34923009.08, 8.528398420000002, 8.562009040000001, 177.15548796, 3.643136682829724, 72.66057681332904, 0.3757955990655879, 13297687.19392669, 3.895599843601784, 3.835131746363534, 90.02768519854139, 0.15),
    @(20, 1, 33.676, 650.096, 1.33198666, 51082170.264, 5.079517299999999, 8.388601679999999, 285.70748816, 6.078232788208701, 113.2815793341215, 0.2257885204055739, 17000659.66582162, 2.146417392504683, 3.483317270665682, 138.6567020361459, 0.85)
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}
